# Auto-generated Excel COM-interop script applying the Seraph_Profits market-data refresh diff.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 548.44446
$ws.Range("I19").Value = 517.7273
$ws.Range("J19").Value = 596.7143
$ws.Range("K19").Value = 517.7273
$ws.Range("L19").Value = 596.7143
$ws.Range("M19").Value = -342.7273
$ws.Range("N19").Value = -946.7143
$ws.Range("H40").Value = 2178.5
$ws.Range("I40").Value = 1928.4286
$ws.Range("J40").Value = 2428.5715
$ws.Range("K40").Value = 1928.4286
$ws.Range("L40").Value = 2428.5715
$ws.Range("M40").Value = -1753.4286
$ws.Range("N40").Value = -2778.5715
$ws.Range("H70").Value = 4956.7334
$ws.Range("I70").Value = 1666
$ws.Range("K70").Value = 4998
$ws.Range("M70").Value = -4728
$ws.Range("H73").Value = 4956.7334
$ws.Range("I73").Value = 1666
$ws.Range("K73").Value = 4998
$ws.Range("M73").Value = -4062
$ws.Range("H76").Value = 4628.4287
$ws.Range("I76").Value = 3166.6667
$ws.Range("J76").Value = 5724.75
$ws.Range("K76").Value = 3166.6667
$ws.Range("L76").Value = 5724.75
$ws.Range("M76").Value = -2851.6667
$ws.Range("N76").Value = -6354.75
$ws.Range("H79").Value = 4628.4287
$ws.Range("I79").Value = 3166.6667
$ws.Range("J79").Value = 5724.75
$ws.Range("K79").Value = 3166.6667
$ws.Range("L79").Value = 5724.75
$ws.Range("M79").Value = -2074.6667
$ws.Range("N79").Value = -7908.75
$ws.Range("H132").Value = 2124.1765
$ws.Range("I132").Value = 2340.7334
$ws.Range("K132").Value = 7022.2002
$ws.Range("M132").Value = -4492.2002
$ws.Range("H137").Value = 2382.6086
$ws.Range("I137").Value = 1771.5714
$ws.Range("J137").Value = 3333.111
$ws.Range("K137").Value = 5314.7142
$ws.Range("L137").Value = 9999.332999999999
$ws.Range("M137").Value = -2764.7142
$ws.Range("N137").Value = -15099.333
$ws.Range("H138").Value = 4783.269
$ws.Range("I138").Value = 1127.8
$ws.Range("K138").Value = 3383.4
$ws.Range("M138").Value = 1756.6
$ws.Range("H141").Value = 2192.158
$ws.Range("I141").Value = 1615.125
$ws.Range("K141").Value = 4845.375
$ws.Range("M141").Value = 334.625

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 30000000
$ws.Range("I13").Value = 30000000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 30000000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -29999856
$ws.Range("H46").Value = 15106
$ws.Range("J46").Value = 13030.25
$ws.Range("L46").Value = 13030.25
$ws.Range("N46").Value = -13668.25
$ws.Range("H61").Value = 2662.7144
$ws.Range("I61").Value = 2528.2
$ws.Range("K61").Value = 2528.2
$ws.Range("M61").Value = -2316.2
$ws.Range("H74").Value = 1347.2285
$ws.Range("I74").Value = 894.0645
$ws.Range("K74").Value = 894.0645
$ws.Range("M74").Value = -20.06449999999995
$ws.Range("H77").Value = 1347.2285
$ws.Range("I77").Value = 894.0645
$ws.Range("K77").Value = 4470.3225
$ws.Range("M77").Value = -102.3225000000002
$ws.Range("H132").Value = 1683.95
$ws.Range("I132").Value = 1584.0588
$ws.Range("K132").Value = 4752.1764
$ws.Range("M132").Value = -2222.1764
$ws.Range("H136").Value = 2662.7144
$ws.Range("I136").Value = 2528.2
$ws.Range("K136").Value = 7584.599999999999
$ws.Range("M136").Value = -5034.599999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2644.7097
$ws.Range("I134").Value = 2439.0454
$ws.Range("K134").Value = 7317.1362
$ws.Range("M134").Value = -4782.1362

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 339.4
$ws.Range("I2").Value = 339.4
$ws.Range("K2").Value = 339.4
$ws.Range("M2").Value = -226.4
$ws.Range("H58").Value = 2974.4
$ws.Range("I58").Value = 1328.4117
$ws.Range("K58").Value = 1328.4117
$ws.Range("M58").Value = -1125.4117
$ws.Range("H132").Value = 2400.1135
$ws.Range("J132").Value = 5343.857
$ws.Range("L132").Value = 16031.571
$ws.Range("N132").Value = -21091.571
$ws.Range("H134").Value = 2155.6667
$ws.Range("I134").Value = 1746.2354
$ws.Range("J134").Value = 3150
$ws.Range("K134").Value = 5238.706200000001
$ws.Range("L134").Value = 9450
$ws.Range("M134").Value = -2703.706200000001
$ws.Range("N134").Value = -14520
$ws.Range("H136").Value = 2974.4
$ws.Range("I136").Value = 1328.4117
$ws.Range("K136").Value = 3985.2351
$ws.Range("M136").Value = -1435.2351

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 26264026
$ws.Range("I4").Value = 32437882
$ws.Range("K4").Value = 97313646
$ws.Range("M4").Value = -97313534
$ws.Range("H12").Value = 514
$ws.Range("I12").Value = 457.83334
$ws.Range("K12").Value = 1373.50002
$ws.Range("M12").Value = -1200.50002
$ws.Range("H23").Value = 199
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H107").Value = 390.7742
$ws.Range("I107").Value = 250.125
$ws.Range("K107").Value = 750.375
$ws.Range("M107").Value = 1169.625

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2319.1614
$ws.Range("I132").Value = 1805.8948
$ws.Range("J132").Value = 3131.8333
$ws.Range("K132").Value = 5417.6844
$ws.Range("L132").Value = 9395.499899999999
$ws.Range("M132").Value = -2887.6844
$ws.Range("N132").Value = -14455.4999

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2479
$ws.Range("I30").Value = 4458
$ws.Range("K30").Value = 4458
$ws.Range("M30").Value = -4350
$ws.Range("H46").Value = 3055.4443
$ws.Range("I46").Value = 2750
$ws.Range("J46").Value = 3666.3333
$ws.Range("K46").Value = 2750
$ws.Range("L46").Value = 3666.3333
$ws.Range("M46").Value = -2562
$ws.Range("N46").Value = -4042.3333
$ws.Range("H61").Value = 6999.5
$ws.Range("J61").Value = 4000
$ws.Range("L61").Value = 4000
$ws.Range("N61").Value = -4404
$ws.Range("H113").Value = 6999.5
$ws.Range("J113").Value = 4000
$ws.Range("L113").Value = 4000
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 7523.5
$ws.Range("I122").Value = 9996.6
$ws.Range("J122").Value = 5050.4
$ws.Range("K122").Value = 29989.8
$ws.Range("L122").Value = 15151.2
$ws.Range("M122").Value = -27539.8
$ws.Range("N122").Value = -20051.2
$ws.Range("H132").Value = 3345.0908
$ws.Range("I132").Value = 2941.25
$ws.Range("K132").Value = 8823.75
$ws.Range("M132").Value = -6293.75

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("H122").Value = 3615.889
$ws.Range("I122").Value = 4220.4287
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 12661.2861
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -10211.2861
$ws.Range("H132").Value = 53583.668
$ws.Range("I132").Value = 68343.92999999999
$ws.Range("J132").Value = 1922.75
$ws.Range("K132").Value = 205031.79
$ws.Range("L132").Value = 5768.25
$ws.Range("M132").Value = -202501.79
$ws.Range("N132").Value = -10828.25

# ----- Special cases: cells removed entirely (no longer populated) -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N2").ClearContents()

# ----- Special case: new cell populated where none existed before -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N122").Value = -9400
